$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44694
$ws.Range("J2").Value = 480
$ws.Range("K2").Value = 17500
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17750
$ws.Range("P2").Value = 710

$ws.Range("D3").Value = 44858
$ws.Range("J3").Value = 500

$ws.Range("D4").Value = 44883
$ws.Range("J4").Value = 380
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7500
$ws.Range("P4").Value = 300

$ws.Range("D6").Value = 44876
$ws.Range("J6").Value = 460
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6500
$ws.Range("P6").Value = 260

$ws.Range("D7").Value = 44803
$ws.Range("J7").Value = 600

$ws.Range("D8").Value = 44824
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 8500
$ws.Range("P8").Value = 340

$ws.Range("D9").Value = 44848
$ws.Range("J9").Value = 800
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 7500
$ws.Range("P9").Value = 300

$ws.Range("D10").Value = 44721
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 14500
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14750
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 590

$ws.Range("D11").Value = 44816
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 9500
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 9750
$ws.Range("P11").Value = 390

$ws.Range("D12").Value = 44799
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 10500
$ws.Range("P12").Value = 420

$ws.Range("D13").Value = 44798
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 10500
$ws.Range("L13").Value = 11000
$ws.Range("M13").Value = 10750
$ws.Range("P13").Value = 430

$ws.Range("D14").Value = 44756
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("P14").Value = 580

$ws.Range("D15").Value = 44827
$ws.Range("J15").Value = 700

$ws.Range("D16").Value = 44377
$ws.Range("K16").Value = 12500
$ws.Range("L16").Value = 13000
$ws.Range("M16").Value = 12750
$ws.Range("P16").Value = 510

$ws.Range("D17").Value = 44855
$ws.Range("J17").Value = 540
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 7500
$ws.Range("P17").Value = 300

$ws.Range("D18").Value = 44825
$ws.Range("J18").Value = 440
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 8500
$ws.Range("O18").Value = "Provincia del Elquí"
$ws.Range("P18").Value = 340

$ws.Range("D20").Value = 44817
$ws.Range("J20").Value = 440
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 9500
$ws.Range("P20").Value = 380

$ws.Range("D21").Value = 44781
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 11000
$ws.Range("M21").Value = 10500
$ws.Range("P21").Value = 420

$ws.Range("D22").Value = 44873
$ws.Range("J22").Value = 540
$ws.Range("K22").Value = 6000
$ws.Range("L22").Value = 7000
$ws.Range("M22").Value = 6500
$ws.Range("O22").Value = "Provincia del Elquí"
$ws.Range("P22").Value = 260

$ws.Range("D23").Value = 44376
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = 12500
$ws.Range("P23").Value = 500

$ws.Range("D24").Value = 44316
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 17000
$ws.Range("M24").Value = 16500
$ws.Range("P24").Value = 660

$ws.Range("D25").Value = 44386
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 11000
$ws.Range("L25").Value = 12000
$ws.Range("M25").Value = 11500
$ws.Range("P25").Value = 460

$ws.Range("D26").Value = 44473
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 8500
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = 8750
$ws.Range("P26").Value = 350

$ws.Range("D27").Value = 44690
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 17500
$ws.Range("P27").Value = 700

$ws.Range("D28").Value = 44356
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 13000
$ws.Range("L28").Value = 14000
$ws.Range("M28").Value = 13500
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 540

$ws.Range("D29").Value = 44384
$ws.Range("J29").Value = 560
$ws.Range("K29").Value = 11500
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = 11750
$ws.Range("P29").Value = 470

$ws.Range("D30").Value = 44809
$ws.Range("J30").Value = 520
$ws.Range("K30").Value = 9500
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = 9750
$ws.Range("P30").Value = 390

$ws.Range("D31").Value = 44797
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 11000
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = 11500
$ws.Range("P31").Value = 460

$ws.Range("D32").Value = 44874
$ws.Range("J32").Value = 500

$ws.Range("D33").Value = 44881
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 7000
$ws.Range("L33").Value = 8000
$ws.Range("M33").Value = 7500
$ws.Range("P33").Value = 300

$ws.Range("D34").Value = 44847
$ws.Range("J34").Value = 520
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = 7500
$ws.Range("P34").Value = 300

$ws.Range("D35").Value = 44372
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 13000
$ws.Range("L35").Value = 14000
$ws.Range("M35").Value = 13500
$ws.Range("P35").Value = 540

$ws.Range("D36").Value = 44811
$ws.Range("J36").Value = 400
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 10500
$ws.Range("M36").Value = 10250
$ws.Range("P36").Value = 410

$ws.Range("D37").Value = 44425
$ws.Range("J37").Value = 400

$ws.Range("D38").Value = 44446
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 11000
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = 11500
$ws.Range("P38").Value = 460

$ws.Range("D39").Value = 44370
$ws.Range("J39").Value = 520
$ws.Range("K39").Value = 13000
$ws.Range("L39").Value = 14000
$ws.Range("M39").Value = 13500
$ws.Range("O39").Value = "Provincia del Elquí"
$ws.Range("P39").Value = 540

$ws.Range("D40").Value = 44466
$ws.Range("J40").Value = 400
$ws.Range("K40").Value = 9500
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = 9750
$ws.Range("P40").Value = 390

$ws.Range("D42").Value = 44837
$ws.Range("J42").Value = 520
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 9000
$ws.Range("M42").Value = 8500
$ws.Range("P42").Value = 340

$ws.Range("D43").Value = 44714
$ws.Range("J43").Value = 400
$ws.Range("K43").Value = 14000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = 14500
$ws.Range("O43").Value = "Provincia de Limarí"
$ws.Range("P43").Value = 580

